$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new "3keV" absorption-coefficient block next to the existing C.(x) columns (rows 31-33) ---
# Keep the order of first-use below so new shared-string indices line up with the target file
# (K31/L31 -> C.(5)/3keV, L32 -> absorption label, then back up to A13, then B34/B36/C36).
$ws.Range("K31").Value = "C.(5)"
$ws.Range("L31").Value = "3keV"
$ws.Range("M31").Formula = "=4394.3*EXP(-458.6*0.003)"
$ws.Range("N31").Value = "cm^2/g"

$ws.Range("L32").Value = "吸収係数/cm"
$ws.Range("M32").Formula = "=M31*M27"

$ws.Range("L33").Value = "C.(5)"
$ws.Range("M33").Formula = "=1/M32"
$ws.Range("N33").Value = "(cm)"

# --- electron rest mass block used by the new velocity table ---
$ws.Range("A13").Value = "電子質量(MeV)"
$ws.Range("A14").Value = 0.51

# --- new "D(1)" / velocity-vs-kinetic-energy table (rows 34-40) ---
$ws.Range("B34").Value = "D(1) "

$ws.Range("B36").Value = "運動エネルギーT (MeV) "
$ws.Range("C36").Value = "速度"

$ws.Range("B37").Value = 0.5
$ws.Range("C37").Formula = "=SQRT(1-(`$A`$14/(B37+`$A`$14))^2)"

$ws.Range("B38").Value = 1
$ws.Range("C38").Formula = "=SQRT(1-(`$A`$14/(B38+`$A`$14))^2)"

$ws.Range("B39").Value = 1.5
$ws.Range("C39").Formula = "=SQRT(1-(`$A`$14/(B39+`$A`$14))^2)"

$ws.Range("B40").Value = 2
$ws.Range("C40").Formula = "=SQRT(1-(`$A`$14/(B40+`$A`$14))^2)"

# --- column widths for the two newly used columns on the left ---
$ws.Columns.Item(1).ColumnWidth = 27.15
$ws.Columns.Item(2).ColumnWidth = 19.4

# --- view state: scroll so row 12 is at the top and select the new F35 cell, like the author left it ---
$excel.Goto($ws.Range("A12"), $true) | Out-Null
$ws.Range("F35").Select() | Out-Null
